$d = $word.ActiveDocument

# Fill in the name placeholder.
$d.Content.Find.Execute("[NAME]", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Sreeraaghav Raja", 2) | Out-Null

# The original "[NAME]<line-break>[DATE]" lived in a single paragraph
# separated by a manual line break (^l). Turn that manual break into a
# real paragraph break (^p) so the name and date become two paragraphs,
# inheriting the same paragraph/run formatting. Only the first manual
# break (the one in the title block) should be touched, so replace just
# one occurrence instead of every manual break in the document.
$d.Content.Find.Execute("^l", $false, $false, $false, $false, $false, `
    $true, 1, $false, "^p", 1) | Out-Null

# Fill in the date placeholder.
$d.Content.Find.Execute("[DATE]", $false, $false, $false, $false, $false, `
    $true, 1, $false, "1/19/2026", 2) | Out-Null
